$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in A1:B4
$ws.Range("A1").Value = -0.015640352534859642
$ws.Range("B1").Value = 0.01564035172969324

$ws.Range("A2").Value = 0.030706584441672825
$ws.Range("B2").Value = -0.030706585194948031

$ws.Range("A3").Value = -0.067844991244612934
$ws.Range("B3").Value = 0.067844990514494455

$ws.Range("A4").Value = -0.01125446203257224
$ws.Range("B4").Value = 0.011254461268442047

# Add new row 5
$ws.Range("A5").Value = 0.028235971966561822
$ws.Range("B5").Value = -0.028235972757080691

# Update column widths (target stored width 14.42578125 for both columns;
# the COM layer quantizes ColumnWidth to an integer pixel grid, so 13.6 is
# the closest input that rounds to the nearest achievable stored width)
$ws.Columns.Item(1).ColumnWidth = 13.6
$ws.Columns.Item(2).ColumnWidth = 13.6
